$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to remain text so values like "504.53" are not
# auto-converted to numbers by Excel's type inference.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "56.675.66"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "2.386.60"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "504.53"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").Value = "132.85"
$ws.Range("E6").Value = "  +3.21%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").Value = "2.392.13"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").Value = "0.0975"
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("E12").Value = "  +2.71%  "
$ws.Range("D13").Value = "4.70"
$ws.Range("E13").Value = "  +2.10%  "
$ws.Range("D14").Value = "2.810.70"
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").Value = "56.618.73"
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").Value = "21.71"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("E17").Value = "  +1.32%  "
$ws.Range("D18").Value = "2.399.25"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("D19").Value = "10.18"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").Value = "4.05"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").Value = "309.72"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").Value = "5.63"
$ws.Range("E24").Value = "  -3.81%  "
$ws.Range("D25").Value = "66.26"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "0.150"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("B28").Value = "Polygon"
$ws.Range("C28").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D28").Value = "0.369"
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("E29").Value = "  +1.79%  "
$ws.Range("D30").Value = "175.33"
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("E31").Value = "  +2.25%  "
$ws.Range("E32").Value = "  -0.59%  "
$ws.Range("E33").Value = "  +2.33%  "
$ws.Range("E34").Value = "  -3.65%  "
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").Value = "17.78"
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("D39").Value = "3.80"
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("D40").Value = "36.79"
$ws.Range("E40").Value = "  +2.63%  "
$ws.Range("D41").Value = "0.818"
$ws.Range("E41").Value = "  +6.26%  "
$ws.Range("D42").Value = "1.43"
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("D43").Value = "132.46"
$ws.Range("E43").Value = "  +1.76%  "
$ws.Range("E44").Value = "  +1.10%  "
$ws.Range("D45").Value = "4.82"
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").Value = "0.0907"
$ws.Range("D48").Value = "246.74"
$ws.Range("E48").Value = "  -2.70%  "
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").Value = "0.0209"
$ws.Range("E50").Value = "  +1.55%  "
$ws.Range("D51").Value = "17.15"
$ws.Range("E51").Value = "  +7.27%  "

# Restore default (no explicit number format) on the Price column so the
# cell styling matches the original workbook.
$ws.Range("D2:D51").ClearFormats()
